$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.035150047760299
$ws.Range("D2").Value = 1.035217875524962
$ws.Range("E2").Value = 1.038796085417311
$ws.Range("F2").Value = 1.033814328158322
$ws.Range("I2").Value = 1.031332201869268
$ws.Range("J2").Value = 1.040265397003256
$ws.Range("K2").Value = 1.038015111781827
$ws.Range("L2").Value = 1.041583099354447
$ws.Range("M2").Value = 1.036615595094452
$ws.Range("N2").Value = 1.04174269318714
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.037332733767326
$ws.Range("D3").Value = 1.036802516734667
$ws.Range("E3").Value = 1.040929472981604
$ws.Range("F3").Value = 1.036644642044381
$ws.Range("I3").Value = 1.031810258412847
$ws.Range("J3").Value = 1.042085291554179
$ws.Range("K3").Value = 1.039407041643166
$ws.Range("L3").Value = 1.043523084670722
$ws.Range("M3").Value = 1.03924958629779
$ws.Range("N3").Value = 1.04356517219704
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.038737378077017
$ws.Range("D4").Value = 1.037821337636674
$ws.Range("E4").Value = 1.042302654616618
$ws.Range("F4").Value = 1.038467218635102
$ws.Range("I4").Value = 1.032115435835513
$ws.Range("J4").Value = 1.043255161443126
$ws.Range("K4").Value = 1.040300705067488
$ws.Range("L4").Value = 1.044770784345284
$ws.Range("M4").Value = 1.040944959960219
$ws.Range("N4").Value = 1.044736703435415
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.039326090025138
$ws.Range("D5").Value = 1.038248112246327
$ws.Range("E5").Value = 1.042878242182896
$ws.Range("F5").Value = 1.039231377513639
$ws.Range("I5").Value = 1.032242748365044
$ws.Range("J5").Value = 1.043745162588803
$ws.Range("K5").Value = 1.040674751830281
$ws.Range("L5").Value = 1.04529353614442
$ws.Range("M5").Value = 1.041655598619763
$ws.Range("N5").Value = 1.045227400438923
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.03942483300041
$ws.Range("D6").Value = 1.038319680270984
$ws.Range("E6").Value = 1.042974787470591
$ws.Range("F6").Value = 1.039359564548192
$ws.Range("I6").Value = 1.03226406731481
$ws.Range("J6").Value = 1.043827330734342
$ws.Range("K6").Value = 1.040737460029199
$ws.Range("L6").Value = 1.04538120513228
$ws.Range("M6").Value = 1.041774796682551
$ws.Range("N6").Value = 1.045309685272653
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.038745251492647
$ws.Range("D7").Value = 1.037827046222044
$ws.Range("E7").Value = 1.042310352263437
$ws.Range("F7").Value = 1.03847743734704
$ws.Range("I7").Value = 1.032117140843488
$ws.Range("J7").Value = 1.043261715933717
$ws.Range("K7").Value = 1.040305709541402
$ws.Range("L7").Value = 1.044777776327775
$ws.Range("M7").Value = 1.040954463710512
$ws.Range("N7").Value = 1.044743267234134
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.035889316897537
$ws.Range("D8").Value = 1.035754786378509
$ws.Range("E8").Value = 1.039518602262367
$ws.Range("F8").Value = 1.03477271019649
$ws.Range("I8").Value = 1.031494630859608
$ws.Range("J8").Value = 1.040882060498105
$ws.Range("K8").Value = 1.038486989381384
$ws.Range("L8").Value = 1.042240323556056
$ws.Range("M8").Value = 1.037507661749536
$ws.Range("N8").Value = 1.042360232414853
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.030795881022993
$ws.Range("D9").Value = 1.032051689681033
$ws.Range("E9").Value = 1.034541704547588
$ws.Range("F9").Value = 1.028174155688182
$ws.Range("I9").Value = 1.030365341112394
$ws.Range("J9").Value = 1.036627971253785
$ws.Range("K9").Value = 1.035227201379011
$ws.Range("L9").Value = 1.037709046822569
$ws.Range("M9").Value = 1.031362475511802
$ws.Range("N9").Value = 1.038100101875856
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.02735654660937
$ws.Range("D10").Value = 1.029546397040198
$ws.Range("E10").Value = 1.031182502954664
$ws.Range("F10").Value = 1.023723913373111
$ws.Range("I10").Value = 1.029590009960894
$ws.Range("J10").Value = 1.033748632537951
$ws.Range("K10").Value = 1.033015213159288
$ws.Range("L10").Value = 1.034645402644122
$ws.Range("M10").Value = 1.027213954086146
$ws.Range("N10").Value = 1.03521667416878
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.025856285806598
$ws.Range("D11").Value = 1.028452468747356
$ws.Range("E11").Value = 1.029717548819672
$ws.Range("F11").Value = 1.021783880076485
$ws.Range("I11").Value = 1.029248789077376
$ws.Range("J11").Value = 1.032491043162088
$ws.Range("K11").Value = 1.03204777964048
$ws.Range("L11").Value = 1.033308099794569
$ws.Range("M11").Value = 1.025404497799655
$ws.Range("N11").Value = 1.033957298871756
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.025297314212309
$ws.Range("D12").Value = 1.028044726732119
$ws.Range("E12").Value = 1.029171785688024
$ws.Range("F12").Value = 1.021061227833
$ws.Range("I12").Value = 1.029121204161981
$ws.Range("J12").Value = 1.032022246384035
$ws.Range("K12").Value = 1.031686948899015
$ws.Range("L12").Value = 1.032809705755955
$ws.Range("M12").Value = 1.024730341738075
$ws.Range("N12").Value = 1.033487836348506
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.025417293691138
$ws.Range("D13").Value = 1.028132253131444
$ws.Range("E13").Value = 1.02928892761168
$ws.Range("F13").Value = 1.02121633265053
$ws.Range("I13").Value = 1.029148609818681
$ws.Range("J13").Value = 1.032122881316743
$ws.Range("K13").Value = 1.031764416013609
$ws.Range("L13").Value = 1.032916688896037
$ws.Range("M13").Value = 1.02487504415334
$ws.Range("N13").Value = 1.033588614194362
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.025810116200482
$ws.Range("D14").Value = 1.028418793618786
$ws.Range("E14").Value = 1.029672469096759
$ws.Range("F14").Value = 1.021724187433059
$ws.Range("I14").Value = 1.029238260094393
$ws.Range("J14").Value = 1.032452326660083
$ws.Range("K14").Value = 1.032017983722021
$ws.Range("L14").Value = 1.033266936543001
$ws.Range("M14").Value = 1.025348813959779
$ws.Range("N14").Value = 1.033918527387878
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.026051919149145
$ws.Range("D15").Value = 1.028595152895528
$ws.Range("E15").Value = 1.02990856623257
$ws.Range("F15").Value = 1.022036821203773
$ws.Range("I15").Value = 1.029293384806645
$ws.Range("J15").Value = 1.032655085817368
$ws.Range("K15").Value = 1.032174017561595
$ws.Range("L15").Value = 1.033482514200359
$ws.Range("M15").Value = 1.025640446205
$ws.Range("N15").Value = 1.034121574486423
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.027455877474833
$ws.Range("D16").Value = 1.029618802071729
$ws.Range("E16").Value = 1.031279503672759
$ws.Range("F16").Value = 1.023852385530582
$ws.Range("I16").Value = 1.029612538624789
$ws.Range("J16").Value = 1.033831862719233
$ws.Range("K16").Value = 1.033079212599629
$ws.Range("L16").Value = 1.034733924919052
$ws.Range("M16").Value = 1.027333759190414
$ws.Range("N16").Value = 1.035300022546467
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.028333558736651
$ws.Range("D17").Value = 1.030258441787809
$ws.Range("E17").Value = 1.032136636538692
$ws.Range("F17").Value = 1.024987693639429
$ws.Range("I17").Value = 1.029811253537893
$ws.Range("J17").Value = 1.034567094983339
$ws.Range("K17").Value = 1.033644414424352
$ws.Range("L17").Value = 1.035515996273294
$ws.Range("M17").Value = 1.028392366703483
$ws.Range("N17").Value = 1.03603629892472
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.028844437149382
$ws.Range("D18").Value = 1.030630655479618
$ws.Range("E18").Value = 1.03263558796827
$ws.Range("F18").Value = 1.025648645900272
$ws.Range("I18").Value = 1.02992663115435
$ws.Range("J18").Value = 1.03499490309446
$ws.Range("K18").Value = 1.033973160617822
$ws.Range("L18").Value = 1.035971133844314
$ws.Range("M18").Value = 1.029008573722329
$ws.Range("N18").Value = 1.036464714572428
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.029018455621186
$ws.Range("D19").Value = 1.03075742306091
$ws.Range("E19").Value = 1.032805549501362
$ws.Range("F19").Value = 1.025873803188395
$ws.Range("I19").Value = 1.029965882569238
$ws.Range("J19").Value = 1.035140599536894
$ws.Range("K19").Value = 1.034085098517272
$ws.Range("L19").Value = 1.036126150513206
$ws.Range("M19").Value = 1.02921847270025
$ws.Range("N19").Value = 1.036610617920523
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.028239501689751
$ws.Range("D20").Value = 1.03018990551998
$ws.Range("E20").Value = 1.032044778065947
$ws.Range("F20").Value = 1.024866016057169
$ws.Range("I20").Value = 1.029789988172803
$ws.Range("J20").Value = 1.03448831947302
$ws.Range("K20").Value = 1.0335838697072
$ws.Range("L20").Value = 1.03543219445382
$ws.Range("M20").Value = 1.028278919068844
$ws.Range("N20").Value = 1.035957411544141
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.025694487351964
$ws.Range("D21").Value = 1.028334453774701
$ws.Range("E21").Value = 1.02955957066308
$ws.Range("F21").Value = 1.021574693729498
$ws.Range("I21").Value = 1.029211883635649
$ws.Range("J21").Value = 1.032355359758031
$ws.Range("K21").Value = 1.031943355581919
$ws.Range("L21").Value = 1.033163843569123
$ws.Range("M21").Value = 1.02520935755845
$ws.Range("N21").Value = 1.033821422781702
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.024084426470176
$ws.Range("D22").Value = 1.027159689890682
$ws.Range("E22").Value = 1.027987657901589
$ws.Range("F22").Value = 1.019493481011611
$ws.Range("I22").Value = 1.028843537379363
$ws.Range("J22").Value = 1.031004585104833
$ws.Range("K22").Value = 1.030903301292901
$ws.Range("L22").Value = 1.031728010811755
$ws.Range("M22").Value = 1.023267542694568
$ws.Range("N22").Value = 1.032468729873561
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.024938907960144
$ws.Range("D23").Value = 1.027783241619341
$ws.Range("E23").Value = 1.028821863812075
$ws.Range("F23").Value = 1.020597918558073
$ws.Range("I23").Value = 1.029039271105062
$ws.Range("J23").Value = 1.031721591549517
$ws.Range("K23").Value = 1.031455481085675
$ws.Range("L23").Value = 1.032490102411353
$ws.Range("M23").Value = 1.024298084267369
$ws.Range("N23").Value = 1.033186754549638
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.028282005287893
$ws.Range("D24").Value = 1.030220876818503
$ws.Range("E24").Value = 1.032086288042301
$ws.Range("F24").Value = 1.024921000789212
$ws.Range("I24").Value = 1.029799598711202
$ws.Range("J24").Value = 1.034523917947954
$ws.Range("K24").Value = 1.033611230119837
$ws.Range("L24").Value = 1.035470064073577
$ws.Range("M24").Value = 1.028330185063939
$ws.Range("N24").Value = 1.035993060572992
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.032120155098049
$ws.Range("D25").Value = 1.033015329255195
$ws.Range("E25").Value = 1.035835430217612
$ws.Range("F25").Value = 1.029888774263944
$ws.Range("I25").Value = 1.030661196364006
$ws.Range("J25").Value = 1.037735206104094
$ws.Range("K25").Value = 1.036076635760891
$ws.Range("L25").Value = 1.03888784815828
$ws.Range("M25").Value = 1.032959997472561
$ws.Range("N25").Value = 1.039208909126654
